# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage so Excel
# does not auto-convert numeric-looking strings (e.g. "539.40")
# into real numbers (which would drop the trailing zero, etc).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# ---- Rows 2-29: price / volume(1h) updates ----
$ws.Range("D2").Value = '60.266.95'
$ws.Range("E2").Value = '  +3.74%  '
$ws.Range("D3").Value = '3.209.44'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue "D5" '539.40'
$ws.Range("E5").Value = '  +0.57%  '
Set-TextValue "D6" '145.28'
$ws.Range("E6").Value = '  +4.49%  '
Set-TextValue "D7" '1.00'
$ws.Range("E7").Value = '  -0.03%  '
Set-TextValue "D8" '0.530'
$ws.Range("E8").Value = '  +3.20%  '
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("E10").Value = '  +4.39%  '
Set-TextValue "D11" '0.432'
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("D12").Value = '3.759.64'
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("E13").Value = '  -0.99%  '
Set-TextValue "D14" '0.0000175'
$ws.Range("E14").Value = '  +3.86%  '
Set-TextValue "D15" '26.15'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").Value = '60.302.66'
$ws.Range("E16").Value = '  +3.59%  '
$ws.Range("D17").Value = '3.234.09'
$ws.Range("E17").Value = '  +2.86%  '
Set-TextValue "D18" '6.25'
$ws.Range("E18").Value = '  +0.30%  '
Set-TextValue "D19" '13.13'
$ws.Range("E19").Value = '  +1.24%  '
Set-TextValue "D20" '8.36'
$ws.Range("E20").Value = '  +2.31%  '
Set-TextValue "D21" '383.31'
$ws.Range("E21").Value = '  +2.22%  '
Set-TextValue "D22" '0.997'
$ws.Range("E22").Value = '  -0.20%  '
Set-TextValue "D23" '0.530'
$ws.Range("E23").Value = '  +2.94%  '
Set-TextValue "D24" '70.22'
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  +2.08%  '
Set-TextValue "D26" '8.85'
$ws.Range("E26").Value = '  +11.43%  '
Set-TextValue "D27" '1.00'
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  +2.97%  '
Set-TextValue "D29" '1.91'
$ws.Range("E29").Value = '  +0.99%  '

# ---- Rows 30-32: ranking reshuffled (RenderToken moves up; NEAR & EthereumClassic shift down one place) ----
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D30" '6.21'
$ws.Range("E30").Value = '  +0.88%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D31" '5.46'
$ws.Range("E31").Value = '  +6.00%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D32" '22.43'
$ws.Range("E32").Value = '  +3.17%  '

# ---- Rows 33-36: price / volume(1h) updates ----
$ws.Range("E33").Value = '  +5.18%  '
Set-TextValue "D34" '6.61'
$ws.Range("E34").Value = '  +5.61%  '
Set-TextValue "D35" '156.91'
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("E36").Value = '  +0.12%  '

# ---- Rows 37-38: ranking swapped (Maker moves up; EnergySwap moves down) ----
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '2.776.75'
$ws.Range("E37").Value = '  +4.97%  '

$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D38" '25.78'
$ws.Range("E38").Value = '  +0.95%  '

# ---- Rows 39-51: price / volume(1h) updates ----
$ws.Range("E39").Value = '  +4.43%  '
$ws.Range("E40").Value = '  +0.26%  '
Set-TextValue "D41" '4.26'
$ws.Range("E41").Value = '  +0.15%  '
Set-TextValue "D42" '39.83'
$ws.Range("E42").Value = '  +3.64%  '
Set-TextValue "D43" '0.725'
$ws.Range("E43").Value = '  +3.43%  '
Set-TextValue "D44" '0.0285'
$ws.Range("E44").Value = '  +4.06%  '
$ws.Range("D45").Value = '3.252.08'
$ws.Range("E45").Value = '  +2.21%  '
$ws.Range("E46").Value = '  +3.05%  '
$ws.Range("E47").Value = '  +0.19%  '
Set-TextValue "D48" '6.17'
$ws.Range("E48").Value = '  -0.96%  '
Set-TextValue "D49" '0.805'
$ws.Range("E49").Value = '  +7.54%  '
Set-TextValue "D50" '20.76'
$ws.Range("E50").Value = '  +2.50%  '
Set-TextValue "D51" '1.00'
$ws.Range("E51").Value = '  +0.00%  '
